$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.04469381780243
$ws.Range("D2").Value = 1.049562856752148
$ws.Range("E2").Value = 1.052251187142056
$ws.Range("F2").Value = 1.062399403692952
$ws.Range("I2").Value = 1.03921230573419
$ws.Range("J2").Value = 1.049757533922906
$ws.Range("K2").Value = 1.052319567876396
$ws.Range("L2").Value = 1.055000436907291
$ws.Range("M2").Value = 1.065120855938856
$ws.Range("N2").Value = 1.020395518248822

# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.045800209629632
$ws.Range("D3").Value = 1.050386816192824
$ws.Range("E3").Value = 1.053204524965673
$ws.Range("F3").Value = 1.063365350029963
$ws.Range("I3").Value = 1.039412402655908
$ws.Range("J3").Value = 1.050510509046007
$ws.Range("K3").Value = 1.052955768235653
$ws.Range("L3").Value = 1.055766222328775
$ws.Range("M3").Value = 1.065901230677546
$ws.Range("N3").Value = 1.020651041769713

# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.04651645562775
$ws.Range("D4").Value = 1.050920084672396
$ws.Range("E4").Value = 1.053821986471795
$ws.Range("F4").Value = 1.063990876322227
$ws.Range("I4").Value = 1.039540602570949
$ws.Range("J4").Value = 1.050997506761563
$ws.Range("K4").Value = 1.053366894470144
$ws.Range("L4").Value = 1.05626170122833
$ws.Range("M4").Value = 1.066406056649275
$ws.Range("N4").Value = 1.020816169151875

# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.046817646135292
$ws.Range("D5").Value = 1.05114429675048
$ws.Range("E5").Value = 1.054081707525044
$ws.Range("F5").Value = 1.064253965265301
$ws.Range("I5").Value = 1.039594192005898
$ws.Range("J5").Value = 1.051202185798443
$ws.Range("K5").Value = 1.05353960268353
$ws.Range("L5").Value = 1.056469991630204
$ws.Range("M5").Value = 1.066618254140741
$ws.Range("N5").Value = 1.020885537402421

# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.046868222082421
$ws.Range("D6").Value = 1.05118194444942
$ws.Range("E6").Value = 1.05412532403949
$ws.Range("F6").Value = 1.06429814593827
$ws.Range("I6").Value = 1.039603171966596
$ws.Range("J6").Value = 1.051236549086177
$ws.Range("K6").Value = 1.05356859355921
$ws.Range("L6").Value = 1.056504963969086
$ws.Range("M6").Value = 1.06665388117928
$ws.Range("N6").Value = 1.020897181622592

# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.04652047983282
$ws.Range("D7").Value = 1.050923080503854
$ws.Range("E7").Value = 1.053825456326908
$ws.Range("F7").Value = 1.063994391267965
$ws.Range("I7").Value = 1.039541319837584
$ws.Range("J7").Value = 1.051000241908619
$ws.Range("K7").Value = 1.053369202714221
$ws.Range("L7").Value = 1.056264484450496
$ws.Range("M7").Value = 1.066408892166395
$ws.Range("N7").Value = 1.02081709625567

# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.045067658502485
$ws.Range("D8").Value = 1.049841294161217
$ws.Range("E8").Value = 1.052573249988216
$ws.Range("F8").Value = 1.062725747000645
$ws.Range("I8").Value = 1.039280193490786
$ws.Range("J8").Value = 1.050012052306338
$ws.Range("K8").Value = 1.05253468586671
$ws.Range("L8").Value = 1.055259244759628
$ws.Range("M8").Value = 1.065384613549547
$ws.Range("N8").Value = 1.020481917685066

# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.04251017276005
$ws.Range("D9").Value = 1.047935934566321
$ws.Range("E9").Value = 1.050371232420074
$ws.Range("F9").Value = 1.060494049934384
$ws.Range("I9").Value = 1.038810294296219
$ws.Range("J9").Value = 1.048269007554874
$ws.Range("K9").Value = 1.051060061170726
$ws.Range("L9").Value = 1.053487631222916
$ws.Range("M9").Value = 1.063578742084334
$ws.Range("N9").Value = 1.019889664695349

# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.040806899423646
$ws.Range("D10").Value = 1.046666331009071
$ws.Range("E10").Value = 1.048906298337511
$ws.Range("F10").Value = 1.059008855527734
$ws.Range("I10").Value = 1.038490481740278
$ws.Range("J10").Value = 1.047105832874099
$ws.Range("K10").Value = 1.050074247959006
$ws.Range("L10").Value = 1.052306409247846
$ws.Range("M10").Value = 1.062374211397703
$ws.Range("N10").Value = 1.019493747974368

# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.040069767083497
$ws.Range("D11").Value = 1.046116737487392
$ws.Range("E11").Value = 1.048272701750249
$ws.Range("F11").Value = 1.058366374824787
$ws.Range("I11").Value = 1.038350449860575
$ws.Range("J11").Value = 1.04660189683479
$ws.Range("K11").Value = 1.049646737380931
$ws.Range("L11").Value = 1.051794896775536
$ws.Range("M11").Value = 1.06185249588353
$ws.Range("N11").Value = 1.019322057379927

# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.039796022437109
$ws.Range("D12").Value = 1.045912617570161
$ws.Range("E12").Value = 1.048037465785076
$ws.Range("F12").Value = 1.058127822288544
$ws.Range("I12").Value = 1.038298202987858
$ws.Range("J12").Value = 1.046414671584177
$ws.Range("K12").Value = 1.049487844047309
$ws.Range("L12").Value = 1.051604893202395
$ws.Range("M12").Value = 1.061658685999501
$ws.Range("N12").Value = 1.019258245530246

# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.039854738923942
$ws.Range("D13").Value = 1.045956400908765
$ws.Range("E13").Value = 1.048087919712934
$ws.Range("F13").Value = 1.058178988382398
$ws.Range("I13").Value = 1.038309420650574
$ws.Range("J13").Value = 1.046454833887339
$ws.Range("K13").Value = 1.049521931585059
$ws.Range("L13").Value = 1.0516456498386
$ws.Range("M13").Value = 1.061700259844773
$ws.Range("N13").Value = 1.019271935120215

# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.040047138061823
$ws.Range("D14").Value = 1.04609986438688
$ws.Range("E14").Value = 1.048253254825252
$ws.Range("F14").Value = 1.058346654089348
$ws.Range("I14").Value = 1.038346135865545
$ws.Range("J14").Value = 1.046586421571989
$ws.Range("K14").Value = 1.049633605182629
$ws.Range("L14").Value = 1.051779191122642
$ws.Range("M14").Value = 1.061836475934656
$ws.Range("N14").Value = 1.019316783453348

# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.040165689343266
$ws.Range("D15").Value = 1.04618826011328
$ws.Range("E15").Value = 1.048355137835566
$ws.Range("F15").Value = 1.058449970857307
$ws.Range("I15").Value = 1.038368726479673
$ws.Range("J15").Value = 1.046667491663518
$ws.Range("K15").Value = 1.049702398145374
$ws.Range("L15").Value = 1.051861469652196
$ws.Range("M15").Value = 1.061920400324342
$ws.Range("N15").Value = 1.019344410918471

# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.040855828783131
$ws.Range("D16").Value = 1.046702808997514
$ws.Range("E16").Value = 1.048948363462098
$ws.Range("F16").Value = 1.05905150792628
$ws.Range("I16").Value = 1.038499742536683
$ws.Range("J16").Value = 1.047139271695566
$ws.Range("K16").Value = 1.050102606825476
$ws.Range("L16").Value = 1.052340355920369
$ws.Range("M16").Value = 1.062408832899941
$ws.Range("N16").Value = 1.019505137131888

# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.041288840820388
$ws.Range("D17").Value = 1.047025613337089
$ws.Range("E17").Value = 1.049320673834542
$ws.Range("F17").Value = 1.059429002130333
$ws.Range("I17").Value = 1.038581510368751
$ws.Range("J17").Value = 1.047435133671117
$ws.Range("K17").Value = 1.050353474180594
$ws.Range("L17").Value = 1.052640739298734
$ws.Range("M17").Value = 1.062715174995836
$ws.Range("N17").Value = 1.01960588799488

# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.041541447675608
$ws.Range("D18").Value = 1.047213914389931
$ws.Range("E18").Value = 1.049537906537845
$ws.Range("F18").Value = 1.059649247733016
$ws.Range("I18").Value = 1.038629054504977
$ws.Range("J18").Value = 1.047607678534632
$ws.Range("K18").Value = 1.050499738421023
$ws.Range("L18").Value = 1.052815944391388
$ws.Range("M18").Value = 1.06289384506588
$ws.Range("N18").Value = 1.019664629554415

# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.041627586587284
$ws.Range("D19").Value = 1.047278122666091
$ws.Range("E19").Value = 1.049611989197927
$ws.Range("F19").Value = 1.059724355888845
$ws.Range("I19").Value = 1.038645240435946
$ws.Range("J19").Value = 1.047666507387894
$ws.Range("K19").Value = 1.050549600145763
$ws.Range("L19").Value = 1.052875684206708
$ws.Range("M19").Value = 1.06295476454052
$ws.Range("N19").Value = 1.019684654711692

# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.041242378772184
$ws.Range("D20").Value = 1.046990977957306
$ws.Range("E20").Value = 1.049280721174142
$ws.Range("F20").Value = 1.059388494404217
$ws.Range("I20").Value = 1.038572752932357
$ws.Range("J20").Value = 1.047403393212807
$ws.Range("K20").Value = 1.050326564943961
$ws.Range("L20").Value = 1.052608511359792
$ws.Range("M20").Value = 1.062682308842812
$ws.Range("N20").Value = 1.019595080935988

# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.039990479682229
$ws.Range("D21").Value = 1.046057617312594
$ws.Range("E21").Value = 1.048204564715177
$ws.Range("F21").Value = 1.058297278135087
$ws.Range("I21").Value = 1.038335330574482
$ws.Range("J21").Value = 1.046547673408244
$ws.Range("K21").Value = 1.04960072276299
$ws.Range("L21").Value = 1.051739866680531
$ws.Range("M21").Value = 1.06179636428201
$ws.Range("N21").Value = 1.019303577793634

# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.039203702920934
$ws.Range("D22").Value = 1.045470913182088
$ws.Range("E22").Value = 1.047528579173681
$ws.Range("F22").Value = 1.057611727814112
$ws.Range("I22").Value = 1.038184706925416
$ws.Range("J22").Value = 1.04600941145499
$ws.Range("K22").Value = 1.049143796677779
$ws.Range("L22").Value = 1.051193686122882
$ws.Range("M22").Value = 1.061239211500996
$ws.Range("N22").Value = 1.019120076573661

# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.039620755896489
$ws.Range("D23").Value = 1.045781922989371
$ws.Range("E23").Value = 1.047886871405794
$ws.Range("F23").Value = 1.057975099590351
$ws.Range("I23").Value = 1.038264682960489
$ws.Range("J23").Value = 1.046294776708023
$ws.Range("K23").Value = 1.04938607484117
$ws.Range("L23").Value = 1.051483229442947
$ws.Range("M23").Value = 1.061534580381477
$ws.Range("N23").Value = 1.019217375000076

# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.0412633728534
$ws.Range("D24").Value = 1.047006628148996
$ws.Range("E24").Value = 1.049298773844178
$ws.Range("F24").Value = 1.059406797917761
$ws.Range("I24").Value = 1.038576710503302
$ws.Range("J24").Value = 1.047417735441217
$ws.Range("K24").Value = 1.050338724262593
$ws.Range("L24").Value = 1.052623073788761
$ws.Range("M24").Value = 1.062697159686149
$ws.Range("N24").Value = 1.019599964257066

# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.043171040639999
$ws.Range("D25").Value = 1.048428406515655
$ws.Range("E25").Value = 1.050939966628265
$ws.Range("F25").Value = 1.061070540980312
$ws.Range("I25").Value = 1.038932929562783
$ws.Range("J25").Value = 1.048719829456808
$ws.Range("K25").Value = 1.051441770109711
$ws.Range("L25").Value = 1.053945663305369
$ws.Range("M25").Value = 1.064045713997332
$ws.Range("N25").Value = 1.020042967734641
